# Refresh cryptocurrency Price (D) and Volume(1h) (E) columns with
# the latest values from the upstream feed, keeping Price cells as text
# (leading apostrophe forces Excel to store the value as a string instead
# of auto-converting numeric-looking text such as "609.40" into a number).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.225.10"
$ws.Range("E2").Value = "  +0.40%  "
$ws.Range("D3").Value = "'3.182.65"
$ws.Range("E3").Value = "  -0.92%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "'609.40"
$ws.Range("E5").Value = "  +1.05%  "
$ws.Range("D6").Value = "'155.34"
$ws.Range("E6").Value = "  +1.07%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "'3.180.41"
$ws.Range("E8").Value = "  -0.97%  "
$ws.Range("D9").Value = "'0.546"
$ws.Range("E9").Value = "  +2.36%  "
$ws.Range("D10").Value = "'0.159"
$ws.Range("E10").Value = "  -0.84%  "
$ws.Range("D11").Value = "'5.68"
$ws.Range("E11").Value = "  -7.09%  "
$ws.Range("D12").Value = "'0.519"
$ws.Range("E12").Value = "  +1.56%  "
$ws.Range("E13").Value = "  -1.02%  "
$ws.Range("D14").Value = "'38.50"
$ws.Range("E14").Value = "  -2.30%  "
$ws.Range("D15").Value = "'3.698.24"
$ws.Range("E15").Value = "  -1.16%  "
$ws.Range("D16").Value = "'66.223.14"
$ws.Range("E16").Value = "  +0.18%  "
$ws.Range("D17").Value = "'7.45"
$ws.Range("E17").Value = "  -0.69%  "
$ws.Range("D18").Value = "'3.175.00"
$ws.Range("E18").Value = "  -1.29%  "
$ws.Range("E19").Value = "  +1.12%  "
$ws.Range("D20").Value = "'511.49"
$ws.Range("E20").Value = "  +0.17%  "
$ws.Range("D21").Value = "'15.44"
$ws.Range("E21").Value = "  -0.25%  "
$ws.Range("D22").Value = "'0.732"
$ws.Range("E22").Value = "  -0.90%  "
$ws.Range("D23").Value = "'8.04"
$ws.Range("E23").Value = "  -0.52%  "
$ws.Range("D24").Value = "'14.86"
$ws.Range("E24").Value = "  -4.24%  "
$ws.Range("D25").Value = "'84.74"
$ws.Range("E25").Value = "  -0.42%  "
$ws.Range("E26").Value = "  +0.23%  "
$ws.Range("D27").Value = "'3.02"
$ws.Range("E27").Value = "  -0.05%  "
$ws.Range("D28").Value = "'9.19"
$ws.Range("E28").Value = "  -0.85%  "
$ws.Range("D29").Value = "'2.39"
$ws.Range("E29").Value = "  +4.54%  "
$ws.Range("D30").Value = "'7.22"
$ws.Range("E30").Value = "  +5.92%  "
$ws.Range("D31").Value = "'3.01"
$ws.Range("E31").Value = "  +5.06%  "
$ws.Range("D32").Value = "'28.05"
$ws.Range("E32").Value = "  -0.46%  "
$ws.Range("E33").Value = "  -0.01%  "
$ws.Range("E34").Value = "  -1.64%  "
$ws.Range("D35").Value = "'6.54"
$ws.Range("E35").Value = "  -0.97%  "
$ws.Range("D36").Value = "'506.39"
$ws.Range("E36").Value = "  +4.57%  "
$ws.Range("D37").Value = "'55.12"
$ws.Range("E37").Value = "  -0.12%  "
$ws.Range("E38").Value = "  -2.54%  "
$ws.Range("D39").Value = "'0.0422"
$ws.Range("E39").Value = "  +0.50%  "
$ws.Range("E40").Value = "  +6.39%  "
$ws.Range("D41").Value = "'8.81"
$ws.Range("E41").Value = "  -1.41%  "
$ws.Range("D42").Value = "'0.0₃0686"
$ws.Range("E42").Value = "  +7.19%  "
$ws.Range("E43").Value = "  -3.63%  "
$ws.Range("E44").Value = "  -2.18%  "
$ws.Range("D45").Value = "'2.44"
$ws.Range("E45").Value = "  -0.82%  "
$ws.Range("D46").Value = "'2.828.41"
$ws.Range("E46").Value = "  -4.03%  "
$ws.Range("D47").Value = "'28.16"
$ws.Range("E47").Value = "  -1.85%  "
$ws.Range("E48").Value = "  +2.40%  "
$ws.Range("E50").Value = "  +0.65%  "
$ws.Range("E51").Value = "  +6.34%  "
